# Swap the Batch (B), MRP (E), Qty (F) and Value (G) columns between pairs
# of adjacent rows that represent two different batches of the same item.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(229, 230),
    @(364, 365),
    @(366, 367),
    @(372, 373),
    @(375, 376),
    @(380, 381),
    @(382, 383),
    @(385, 386),
    @(442, 443),
    @(463, 464),
    @(572, 573)
)

$cols = @("B", "E", "F", "G")

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"

        $val1 = $ws.Range($addr1).Value2
        $val2 = $ws.Range($addr2).Value2

        $ws.Range($addr1).Value2 = $val2
        $ws.Range($addr2).Value2 = $val1
    }
}
